# =========================================================================
# Update "广州-漫展信息.xlsx" per commit "Update gh-pages to output
# generated at 456a3b4"
# =========================================================================

$wb = $excel.ActiveWorkbook

function Set-NewRowFormat {
    param($ws, [int]$row)
    # The Insert() call sometimes leaves the new row's first cell with a
    # slightly different style than its neighbours; copy the format from
    # the row above so it matches the rest of the sequential-index column.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)  # xlPasteFormats
}

# -------------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple value-only updates on rows that are not affected by the
# forthcoming row insertion (rows 1-21).
$ws1.Range("F2").Value  = 42
$ws1.Range("F5").Value  = 981

$ws1.Range("C6").Value  = "广州·运动番ONLY（取消）"
$ws1.Range("F6").Value  = 371
$ws1.Range("G6").Value  = "不可售"

$ws1.Range("F9").Value  = 1465
$ws1.Range("F10").Value = 132
$ws1.Range("F11").Value = 1358
$ws1.Range("F12").Value = 3014
$ws1.Range("F13").Value = 444
$ws1.Range("F14").Value = 1639
$ws1.Range("F16").Value = 801
$ws1.Range("F17").Value = 242
$ws1.Range("F18").Value = 1399
$ws1.Range("I18").Value = "//i1.hdslb.com/bfs/openplatform/202404/fnUmp06X1713859959456.png"
$ws1.Range("F19").Value = 269

# Insert a brand-new event row at position 22; existing rows 22-27 shift
# down to 23-28.
$ws1.Rows.Item(22).Insert()

$ws1.Range("B22").Value = "2024-06-01"
$ws1.Range("C22").Value = "广州·WIO JUMPONLY3.0"
$ws1.Range("D22").Value = "黄边三横路一街1号 设计殿堂"
$ws1.Range("E22").Value = "2024.06.01 10:00-06.02 18:00"
$ws1.Range("F22").Value = 2
$ws1.Range("G22").Value = 70
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=84722"
$ws1.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202404/FhaZLO921713774163735.jpeg"

Set-NewRowFormat -ws $ws1 -row 22

# Restore the sequential index column (A) for the new row and everything
# that shifted below it.
for ($r = 22; $r -le 28; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# A few of the shifted-down rows also had their "want to go" counter
# ticked up.
$ws1.Range("F24").Value = 6
$ws1.Range("F25").Value = 3508
$ws1.Range("F26").Value = 695
$ws1.Range("F28").Value = 1551

# -------------------------------------------------------------------------
# Sheet "演出" (performances)
# -------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F5").Value  = 29
$ws2.Range("F6").Value  = 11
$ws2.Range("F8").Value  = 24
$ws2.Range("F11").Value = 26

# -------------------------------------------------------------------------
# Sheet "全部类型" (all types) - mirrors the "展览" sheet changes at the
# appropriate (different) row offsets, plus the "演出" sheet changes.
# -------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 42

$ws4.Range("F9").Value  = 29
$ws4.Range("F10").Value = 11
$ws4.Range("F13").Value = 24

$ws4.Range("F15").Value = 981

$ws4.Range("C16").Value = "广州·运动番ONLY（取消）"
$ws4.Range("F16").Value = 371
$ws4.Range("G16").Value = "不可售"

$ws4.Range("F19").Value = 1465
$ws4.Range("F20").Value = 132
$ws4.Range("F21").Value = 1358
$ws4.Range("F22").Value = 3014
$ws4.Range("F23").Value = 444
$ws4.Range("F24").Value = 1639
$ws4.Range("F26").Value = 801
$ws4.Range("F27").Value = 242
$ws4.Range("F28").Value = 1399
$ws4.Range("I28").Value = "//i1.hdslb.com/bfs/openplatform/202404/fnUmp06X1713859959456.png"
$ws4.Range("F29").Value = 269

$ws4.Range("F32").Value = 26

# Insert the same brand-new event row (position 34 on this sheet); rows
# 34-42 shift down to 35-43.
$ws4.Rows.Item(34).Insert()

$ws4.Range("B34").Value = "2024-06-01"
$ws4.Range("C34").Value = "广州·WIO JUMPONLY3.0"
$ws4.Range("D34").Value = "黄边三横路一街1号 设计殿堂"
$ws4.Range("E34").Value = "2024.06.01 10:00-06.02 18:00"
$ws4.Range("F34").Value = 2
$ws4.Range("G34").Value = 70
$ws4.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=84722"
$ws4.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202404/FhaZLO921713774163735.jpeg"

Set-NewRowFormat -ws $ws4 -row 34

for ($r = 34; $r -le 43; $r++) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
}

$ws4.Range("F37").Value = 6
$ws4.Range("F38").Value = 3508
$ws4.Range("F39").Value = 695
$ws4.Range("F41").Value = 1551

Write-Host "All edits applied"
